$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) from column Q (rows 2-6) into the new
# column R so the new cells pick up the same cell styles (borders,
# number formats, fonts) as their left-hand neighbours.
$ws.Range("Q2:Q6").Copy()
$ws.Range("R2:R6").PasteSpecial(-4122)

# Populate the new column's data (2021 figures).
$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 233306
$ws.Range("R5").Value = 3.5
$ws.Range("R6").Value = 30.8

# Match the saved selection/active cell recorded in the sheet view.
$ws.Range("Q15").Select()
